$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

$ws.Range("D2").Value = "42.614.50"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.287.77"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue $ws.Range("D5") "304.81"
$ws.Range("E5").Value = "  +1.62%  "
Set-TextValue $ws.Range("D6") "95.44"
$ws.Range("E6").Value = "  -2.35%  "
$ws.Range("E7").Value = "  -3.48%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -3.38%  "
Set-TextValue $ws.Range("D10") "34.92"
$ws.Range("E10").Value = "  -3.07%  "
Set-TextValue $ws.Range("D11") "0.0782"
$ws.Range("E11").Value = "  -0.74%  "
Set-TextValue $ws.Range("D12") "18.18"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "2.641.81"
$ws.Range("D16").Value = "2.280.32"
$ws.Range("E16").Value = "  -0.63%  "
Set-TextValue $ws.Range("D17") "0.773"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "42.505.55"
$ws.Range("E18").Value = "  -0.97%  "
Set-TextValue $ws.Range("D19") "12.82"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("E21").Value = "  -2.11%  "
Set-TextValue $ws.Range("D22") "66.76"
$ws.Range("E22").Value = "  -3.06%  "
Set-TextValue $ws.Range("D23") "235.42"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  +0.09%  "
Set-TextValue $ws.Range("D28") "24.94"
$ws.Range("E28").Value = "  +0.15%  "
Set-TextValue $ws.Range("D29") "165.66"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  +0.64%  "
Set-TextValue $ws.Range("D31") "8.97"
$ws.Range("E31").Value = "  -1.04%  "
Set-TextValue $ws.Range("D32") "32.29"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  -2.85%  "
Set-TextValue $ws.Range("D36") "17.51"
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D43").Value = "1.991.85"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("E44").Value = "  -3.03%  "
Set-TextValue $ws.Range("D45") "9.97"
$ws.Range("E45").Value = "  -2.43%  "
Set-TextValue $ws.Range("D46") "17.80"
$ws.Range("E46").Value = "  +2.22%  "
Set-TextValue $ws.Range("D47") "2.01"
$ws.Range("E47").Value = "  -9.88%  "
Set-TextValue $ws.Range("D48") "2.76"
$ws.Range("E48").Value = "  -2.17%  "
Set-TextValue $ws.Range("D49") "2.89"
$ws.Range("E49").Value = "  +8.99%  "
Set-TextValue $ws.Range("D50") "53.27"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").Value = "2.507.15"
$ws.Range("E51").Value = "  -0.48%  "
